$wb = $excel.ActiveWorkbook

# Sheet ALC, row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 76.59999999999999  # H33: 68.833336 -> 76.59999999999999
$ws.Cells.Item(33, 9).Value = 80.75  # I33: 73.25 -> 80.75
$ws.Cells.Item(33, 11).Value = 80.75  # K33: 73.25 -> 80.75
$ws.Cells.Item(33, 13).Value = 148.25  # M33: 155.75 -> 148.25

# Sheet ALC, row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 9499  # H64: 9499.25 -> 9499
$ws.Cells.Item(64, 9).Value = 8998  # I64: 8998.5 -> 8998
$ws.Cells.Item(64, 11).Value = 8998  # K64: 8998.5 -> 8998
$ws.Cells.Item(64, 13).Value = -8750  # M64: -8750.5 -> -8750

# Sheet ALC, row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 9499  # H67: 9499.25 -> 9499
$ws.Cells.Item(67, 9).Value = 8998  # I67: 8998.5 -> 8998
$ws.Cells.Item(67, 11).Value = 8998  # K67: 8998.5 -> 8998
$ws.Cells.Item(67, 13).Value = -8140  # M67: -8140.5 -> -8140

# Sheet ALC, row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 5982.4287  # H70: 2575.4 -> 5982.4287
$ws.Cells.Item(70, 9).Value = 1878  # I70: 1439 -> 1878
$ws.Cells.Item(70, 10).Value = 6666.5  # J70: 3333 -> 6666.5
$ws.Cells.Item(70, 11).Value = 5634  # K70: 4317 -> 5634
$ws.Cells.Item(70, 12).Value = 19999.5  # L70: 9999 -> 19999.5
$ws.Cells.Item(70, 13).Value = -5364  # M70: -4047 -> -5364
$ws.Cells.Item(70, 14).Value = -20539.5  # N70: -10539 -> -20539.5

# Sheet ALC, row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 5982.4287  # H73: 2575.4 -> 5982.4287
$ws.Cells.Item(73, 9).Value = 1878  # I73: 1439 -> 1878
$ws.Cells.Item(73, 10).Value = 6666.5  # J73: 3333 -> 6666.5
$ws.Cells.Item(73, 11).Value = 5634  # K73: 4317 -> 5634
$ws.Cells.Item(73, 12).Value = 19999.5  # L73: 9999 -> 19999.5
$ws.Cells.Item(73, 13).Value = -4698  # M73: -3381 -> -4698
$ws.Cells.Item(73, 14).Value = -21871.5  # N73: -11871 -> -21871.5

# Sheet ALC, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 171  # H107: 228.85715 -> 171
$ws.Cells.Item(107, 10).Value = 296  # J107: 498.5 -> 296
$ws.Cells.Item(107, 12).Value = 296  # L107: 498.5 -> 296
$ws.Cells.Item(107, 14).Value = -4136  # N107: -4338.5 -> -4136

# Sheet ALC, row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 3120.8462  # H125: 3126.4285 -> 3120.8462
$ws.Cells.Item(125, 9).Value = 2964.25  # I125: 2982.3076 -> 2964.25
$ws.Cells.Item(125, 11).Value = 26678.25  # K125: 26840.7684 -> 26678.25
$ws.Cells.Item(125, 13).Value = -24218.25  # M125: -24380.7684 -> -24218.25

# Sheet ALC, row 128
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(128, 8).Value = 68963.336  # H128: 68980 -> 68963.336
$ws.Cells.Item(128, 10).Value = 68963.336  # J128: 68980 -> 68963.336
$ws.Cells.Item(128, 12).Value = 68963.336  # L128: 68980 -> 68963.336
$ws.Cells.Item(128, 14).Value = -78923.336  # N128: -78940 -> -78923.336

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2249.1667  # H137: 2358.9443 -> 2249.1667
$ws.Cells.Item(137, 9).Value = 1662  # I137: 1737.0588 -> 1662
$ws.Cells.Item(137, 10).Value = 2905.4119  # J137: 2915.3684 -> 2905.4119
$ws.Cells.Item(137, 11).Value = 4986  # K137: 5211.1764 -> 4986
$ws.Cells.Item(137, 12).Value = 8716.235700000001  # L137: 8746.1052 -> 8716.235700000001
$ws.Cells.Item(137, 13).Value = -2436  # M137: -2661.1764 -> -2436
$ws.Cells.Item(137, 14).Value = -13816.2357  # N137: -13846.1052 -> -13816.2357

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2261.2  # H138: 2226.875 -> 2261.2
$ws.Cells.Item(138, 9).Value = 2284.2856  # I138: 2232.1667 -> 2284.2856
$ws.Cells.Item(138, 10).Value = 2207.3333  # J138: 2211 -> 2207.3333
$ws.Cells.Item(138, 11).Value = 6852.8568  # K138: 6696.500100000001 -> 6852.8568
$ws.Cells.Item(138, 12).Value = 6621.999899999999  # L138: 6633 -> 6621.999899999999
$ws.Cells.Item(138, 13).Value = -1712.8568  # M138: -1556.500100000001 -> -1712.8568
$ws.Cells.Item(138, 14).Value = -16901.9999  # N138: -16913 -> -16901.9999

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 5287.421  # H74: 5554.5 -> 5287.421
$ws.Cells.Item(74, 9).Value = 2646.9167  # I74: 2843.9092 -> 2646.9167
$ws.Cells.Item(74, 11).Value = 2646.9167  # K74: 2843.9092 -> 2646.9167
$ws.Cells.Item(74, 13).Value = -1772.9167  # M74: -1969.9092 -> -1772.9167

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 5287.421  # H77: 5554.5 -> 5287.421
$ws.Cells.Item(77, 9).Value = 2646.9167  # I77: 2843.9092 -> 2646.9167
$ws.Cells.Item(77, 11).Value = 13234.5835  # K77: 14219.546 -> 13234.5835
$ws.Cells.Item(77, 13).Value = -8866.583500000001  # M77: -9851.546 -> -8866.583500000001

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 4768  # H132: 5796.857 -> 4768
$ws.Cells.Item(132, 9).Value = 802.6667  # I132: 805 -> 802.6667
$ws.Cells.Item(132, 10).Value = 7742  # J132: 7793.6 -> 7742
$ws.Cells.Item(132, 11).Value = 2408.0001  # K132: 2415 -> 2408.0001
$ws.Cells.Item(132, 12).Value = 23226  # L132: 23380.8 -> 23226
$ws.Cells.Item(132, 13).Value = 121.9998999999998  # M132: 115 -> 121.9998999999998
$ws.Cells.Item(132, 14).Value = -28286  # N132: -28440.8 -> -28286

# Sheet BSM, row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 181.5  # H80: 170 -> 181.5
$ws.Cells.Item(80, 10).Value = 196.2  # J80: 178.41667 -> 196.2
$ws.Cells.Item(80, 12).Value = 196.2  # L80: 178.41667 -> 196.2
$ws.Cells.Item(80, 14).Value = -2192.2  # N80: -2174.41667 -> -2192.2

# Sheet BSM, row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 181.5  # H83: 170 -> 181.5
$ws.Cells.Item(83, 10).Value = 196.2  # J83: 178.41667 -> 196.2
$ws.Cells.Item(83, 12).Value = 981  # L83: 892.0833500000001 -> 981
$ws.Cells.Item(83, 14).Value = -10965  # N83: -10876.08335 -> -10965

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2046.9714  # H134: 2315.9678 -> 2046.9714
$ws.Cells.Item(134, 9).Value = 1791.2069  # I134: 1953.6923 -> 1791.2069
$ws.Cells.Item(134, 10).Value = 3283.1667  # J134: 4199.8 -> 3283.1667
$ws.Cells.Item(134, 11).Value = 5373.620699999999  # K134: 5861.0769 -> 5373.620699999999
$ws.Cells.Item(134, 12).Value = 9849.500100000001  # L134: 12599.4 -> 9849.500100000001
$ws.Cells.Item(134, 13).Value = -2838.620699999999  # M134: -3326.0769 -> -2838.620699999999
$ws.Cells.Item(134, 14).Value = -14919.5001  # N134: -17669.4 -> -14919.5001

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1928.619  # H58: 1976.55 -> 1928.619
$ws.Cells.Item(58, 9).Value = 1089.0588  # I58: 1096.5 -> 1089.0588
$ws.Cells.Item(58, 11).Value = 1089.0588  # K58: 1096.5 -> 1089.0588
$ws.Cells.Item(58, 13).Value = -886.0588  # M58: -893.5 -> -886.0588

# Sheet CRP, row 81
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(81, 8).Value = 45786  # H81: 45785.5 -> 45786
$ws.Cells.Item(81, 9).Value = 0  # I81: 44325 -> 0
$ws.Cells.Item(81, 10).Value = 45786  # J81: 46272.332 -> 45786
$ws.Cells.Item(81, 11).Value = 0  # K81: 44325 -> 0
$ws.Cells.Item(81, 12).ClearContents()  # L81: 46272.332 -> (removed)
$ws.Cells.Item(81, 13).Value = 45786  # M81: -43327 -> 45786
$ws.Cells.Item(81, 14).Value = -47782  # N81: -48268.332 -> -47782

# Sheet CRP, row 84
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(84, 8).Value = 45786  # H84: 45785.5 -> 45786
$ws.Cells.Item(84, 9).Value = 0  # I84: 44325 -> 0
$ws.Cells.Item(84, 10).Value = 45786  # J84: 46272.332 -> 45786
$ws.Cells.Item(84, 11).Value = 0  # K84: 132975 -> 0
$ws.Cells.Item(84, 12).ClearContents()  # L84: 138816.996 -> (removed)
$ws.Cells.Item(84, 13).Value = 137358  # M84: -127983 -> 137358
$ws.Cells.Item(84, 14).Value = -147342  # N84: -148800.996 -> -147342

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 4779.875  # H134: 4748.4287 -> 4779.875
$ws.Cells.Item(134, 10).Value = 6083.3335  # J134: 6625 -> 6083.3335
$ws.Cells.Item(134, 12).Value = 18250.0005  # L134: 19875 -> 18250.0005
$ws.Cells.Item(134, 14).Value = -23320.0005  # N134: -24945 -> -23320.0005

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 1928.619  # H136: 1976.55 -> 1928.619
$ws.Cells.Item(136, 9).Value = 1089.0588  # I136: 1096.5 -> 1089.0588
$ws.Cells.Item(136, 11).Value = 3267.1764  # K136: 3289.5 -> 3267.1764
$ws.Cells.Item(136, 13).Value = -717.1764000000003  # M136: -739.5 -> -717.1764000000003

# Sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 78295.46000000001  # H4: 98416.39 -> 78295.46000000001
$ws.Cells.Item(4, 9).Value = 143126.47  # I4: 200234.4 -> 143126.47
$ws.Cells.Item(4, 10).Value = 2659.2778  # J4: 2962 -> 2659.2778
$ws.Cells.Item(4, 11).Value = 429379.41  # K4: 600703.2 -> 429379.41
$ws.Cells.Item(4, 12).Value = 7977.8334  # L4: 8886 -> 7977.8334
$ws.Cells.Item(4, 13).Value = -429267.41  # M4: -600591.2 -> -429267.41
$ws.Cells.Item(4, 14).Value = -8201.8334  # N4: -9110 -> -8201.8334

# Sheet CUL, row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 129.52632  # H12: 128.05 -> 129.52632
$ws.Cells.Item(12, 10).Value = 138.14285  # J12: 135.6 -> 138.14285
$ws.Cells.Item(12, 12).Value = 414.42855  # L12: 406.8 -> 414.42855
$ws.Cells.Item(12, 14).Value = -760.4285500000001  # N12: -752.8 -> -760.4285500000001

# Sheet CUL, row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 9).Value = 100.5  # I33: 86.59999999999999 -> 100.5
$ws.Cells.Item(33, 10).Value = 71  # J33: 91 -> 71
$ws.Cells.Item(33, 11).Value = 603  # K33: 519.5999999999999 -> 603
$ws.Cells.Item(33, 12).Value = 426  # L33: 546 -> 426
$ws.Cells.Item(33, 13).Value = -320  # M33: -236.5999999999999 -> -320
$ws.Cells.Item(33, 14).Value = -992  # N33: -1112 -> -992

# Sheet CUL, row 127
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(127, 8).Value = 3000  # H127: 4000 -> 3000
$ws.Cells.Item(127, 10).Value = 3000  # J127: 4000 -> 3000
$ws.Cells.Item(127, 12).Value = 9000  # L127: 12000 -> 9000
$ws.Cells.Item(127, 14).Value = -18920  # N127: -21920 -> -18920

# Sheet CUL, row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 1631.8  # H129: 1733.875 -> 1631.8
$ws.Cells.Item(129, 9).Value = 888.3333  # I129: 976 -> 888.3333
$ws.Cells.Item(129, 10).Value = 2747  # J129: 2997 -> 2747
$ws.Cells.Item(129, 11).Value = 2664.9999  # K129: 2928 -> 2664.9999
$ws.Cells.Item(129, 12).Value = 8241  # L129: 8991 -> 8241
$ws.Cells.Item(129, 13).Value = 2335.0001  # M129: 2072 -> 2335.0001
$ws.Cells.Item(129, 14).Value = -18241  # N129: -18991 -> -18241

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 47551.523  # H132: 51508.906 -> 47551.523
$ws.Cells.Item(132, 10).Value = 5099.2856  # J132: 4739.4 -> 5099.2856
$ws.Cells.Item(132, 12).Value = 15297.8568  # L132: 14218.2 -> 15297.8568
$ws.Cells.Item(132, 14).Value = -20357.8568  # N132: -19278.2 -> -20357.8568

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 8900.286  # H46: 9500 -> 8900.286
$ws.Cells.Item(46, 9).Value = 7200.5  # I46: 0 -> 7200.5
$ws.Cells.Item(46, 10).Value = 9580.200000000001  # J46: 9500 -> 9580.200000000001
$ws.Cells.Item(46, 11).Value = 7200.5  # K46: 0 -> 7200.5
$ws.Cells.Item(46, 12).Value = 9580.200000000001  # L46: 9500 -> 9580.200000000001
$ws.Cells.Item(46, 13).Value = -7012.5  # M46: None -> -7012.5
$ws.Cells.Item(46, 14).Value = -9956.200000000001  # N46: -9876 -> -9956.200000000001

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 5695.5  # H132: 6040.75 -> 5695.5
$ws.Cells.Item(132, 9).Value = 2803.5  # I132: 3494 -> 2803.5
$ws.Cells.Item(132, 11).Value = 8410.5  # K132: 10482 -> 8410.5
$ws.Cells.Item(132, 13).Value = -5880.5  # M132: -7952 -> -5880.5

# Sheet WVR, row 49
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 3336200  # H49: 5002800 -> 3336200
$ws.Cells.Item(49, 9).Value = 4300  # I49: 5600 -> 4300
$ws.Cells.Item(49, 11).Value = 4300  # K49: 5600 -> 4300
$ws.Cells.Item(49, 13).Value = -4070  # M49: -5370 -> -4070

# Sheet WVR, row 82
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(82, 8).Value = 50000  # H82: 0 -> 50000
$ws.Cells.Item(82, 10).Value = 50000  # J82: 0 -> 50000
$ws.Cells.Item(82, 12).Value = 50000  # L82: 0 -> 50000
$ws.Cells.Item(82, 14).Value = -50766  # N82: None -> -50766

# Sheet WVR, row 85
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(85, 8).Value = 50000  # H85: 0 -> 50000
$ws.Cells.Item(85, 10).Value = 50000  # J85: 0 -> 50000
$ws.Cells.Item(85, 12).Value = 50000  # L85: 0 -> 50000
$ws.Cells.Item(85, 14).Value = -52652  # N85: None -> -52652

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 615.6923  # H107: 619.0769 -> 615.6923
$ws.Cells.Item(107, 9).Value = 530.4  # I107: 531.63635 -> 530.4
$ws.Cells.Item(107, 10).Value = 900  # J107: 1100 -> 900
$ws.Cells.Item(107, 11).Value = 1591.2  # K107: 1594.90905 -> 1591.2
$ws.Cells.Item(107, 12).Value = 2700  # L107: 3300 -> 2700
$ws.Cells.Item(107, 13).Value = 328.8000000000002  # M107: 325.09095 -> 328.8000000000002
$ws.Cells.Item(107, 14).Value = -6540  # N107: -7140 -> -6540

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 806.0769  # H113: 935.63635 -> 806.0769
$ws.Cells.Item(113, 9).Value = 1054.75  # I113: 1192.7142 -> 1054.75
$ws.Cells.Item(113, 10).Value = 408.2  # J113: 485.75 -> 408.2
$ws.Cells.Item(113, 11).Value = 3164.25  # K113: 3578.1426 -> 3164.25
$ws.Cells.Item(113, 12).Value = 1224.6  # L113: 1457.25 -> 1224.6
$ws.Cells.Item(113, 13).Value = -994.25  # M113: -1408.1426 -> -994.25
$ws.Cells.Item(113, 14).Value = -5564.6  # N113: -5797.25 -> -5564.6

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3948.6956  # H136: 4233.095 -> 3948.6956
$ws.Cells.Item(136, 9).Value = 2567.2307  # I136: 2859 -> 2567.2307
$ws.Cells.Item(136, 11).Value = 7701.6921  # K136: 8577 -> 7701.6921
$ws.Cells.Item(136, 13).Value = -5151.6921  # M136: -6027 -> -5151.6921
